# Update the "Förändrad" (Changed) date column (C) for all data rows.
# The workbook stores dates as serial numbers; all rows in C2:C246 had the
# same value (46075, i.e. 2026-02-22) and need to be bumped by one day to
# 46076 (2026-02-23), leaving everything else (formatting, other columns)
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 246) { $lastRow = 246 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46075) {
        $cell.Value2 = 46076
    }
}
